# Update the two-digit / one-digit division drill numbers.
# The worksheet table has 20 rows; rows 1, 5, 9, 13, 17 hold the 5 problems
# per block (columns 1-5), the remaining rows are blank answer rows.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $rng = $cell.Range
    # Exclude the trailing end-of-cell marker so we don't clobber it.
    $rng.End = $rng.End - 1
    $rng.Text = $newText
}

# Row 1
Set-CellText $t 1 1 "45÷7="
Set-CellText $t 1 2 "79÷4="
Set-CellText $t 1 3 "44÷7="
Set-CellText $t 1 4 "51÷4="
Set-CellText $t 1 5 "50÷2="

# Row 5
Set-CellText $t 5 1 "89÷3="
Set-CellText $t 5 2 "19÷4="
Set-CellText $t 5 3 "72÷6="
Set-CellText $t 5 4 "71÷3="
Set-CellText $t 5 5 "82÷6="

# Row 9
Set-CellText $t 9 1 "49÷7="
Set-CellText $t 9 2 "61÷2="
Set-CellText $t 9 3 "80÷8="
Set-CellText $t 9 4 "91÷5="
Set-CellText $t 9 5 "93÷6="

# Row 13
Set-CellText $t 13 1 "80÷5="
Set-CellText $t 13 2 "50÷9="
Set-CellText $t 13 3 "45÷8="
Set-CellText $t 13 4 "36÷7="
Set-CellText $t 13 5 "55÷5="

# Row 17
Set-CellText $t 17 1 "23÷2="
Set-CellText $t 17 2 "51÷4="
Set-CellText $t 17 3 "13÷3="
Set-CellText $t 17 4 "20÷3="
Set-CellText $t 17 5 "94÷9="
